# Update Apoe-Lrp1 LR-pair sheet with new TPM-derived values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 32.50235
$ws.Range("H2").Value = 97.50704999999999
$ws.Range("I2").Value = 0.004318312013857221
$ws.Range("J2").Value = 0.004318312013857221
$ws.Range("M2").Value = 6.305846
$ws.Range("N2").Value = 18.917538
$ws.Range("O2").Value = 0.01356150511917599
$ws.Range("P2").Value = 0.01356150511917599
$ws.Range("Q2").Value = 204.9548137381
$ws.Range("R2").Value = 1844.5933236429
$ws.Range("S2").Value = 0.00005856281048212388
$ws.Range("T2").Value = 0.00005856281048212387
$ws.Range("G3").Value = 32.50235
$ws.Range("H3").Value = 97.50704999999999
$ws.Range("I3").Value = 0.004318312013857221
$ws.Range("J3").Value = 0.004318312013857221
$ws.Range("O3").Value = 0.392557056479861
$ws.Range("P3").Value = 0.3925570564798609
$ws.Range("Q3").Value = 5932.708625286818
$ws.Range("R3").Value = 53394.37762758135
$ws.Range("S3").Value = 0.001695183853121411
$ws.Range("T3").Value = 0.001695183853121411
$ws.Range("G4").Value = 32.50235
$ws.Range("H4").Value = 97.50704999999999
$ws.Range("I4").Value = 0.004318312013857221
$ws.Range("J4").Value = 0.004318312013857221
$ws.Range("M4").Value = 127.396393
$ws.Range("N4").Value = 382.189179
$ws.Range("O4").Value = 0.2739817680029065
$ws.Range("P4").Value = 0.2739817680029065
$ws.Range("Q4").Value = 4140.682154023551
$ws.Range("R4").Value = 37266.13938621195
$ws.Range("S4").Value = 0.001183138760344793
$ws.Range("T4").Value = 0.001183138760344793
$ws.Range("G5").Value = 32.50235
$ws.Range("H5").Value = 97.50704999999999
$ws.Range("I5").Value = 0.004318312013857221
$ws.Range("J5").Value = 0.004318312013857221
$ws.Range("M5").Value = 19.42400133333333
$ws.Range("N5").Value = 58.272004
$ws.Range("O5").Value = 0.04177372766745037
$ws.Range("P5").Value = 0.04177372766745036
$ws.Range("Q5").Value = 631.3256897364666
$ws.Range("R5").Value = 5681.931207628199
$ws.Range("S5").Value = 0.0001803919900499507
$ws.Range("T5").Value = 0.0001803919900499507
$ws.Range("G6").Value = 32.50235
$ws.Range("H6").Value = 97.50704999999999
$ws.Range("I6").Value = 0.004318312013857221
$ws.Range("J6").Value = 0.004318312013857221
$ws.Range("M6").Value = 129.3233566666667
$ws.Range("N6").Value = 387.97007
$ws.Range("O6").Value = 0.2781259427306063
$ws.Range("P6").Value = 0.2781259427306062
$ws.Range("Q6").Value = 4203.313001554834
$ws.Range("R6").Value = 37829.8170139935
$ws.Range("S6").Value = 0.001201034599858943
$ws.Range("T6").Value = 0.001201034599858942
$ws.Range("I7").Value = 0.006762540683959845
$ws.Range("J7").Value = 0.006762540683959845
$ws.Range("M7").Value = 6.305846
$ws.Range("N7").Value = 18.917538
$ws.Range("O7").Value = 0.01356150511917599
$ws.Range("P7").Value = 0.01356150511917599
$ws.Range("Q7").Value = 320.9622792030006
$ws.Range("R7").Value = 2888.660512827006
$ws.Range("S7").Value = 0.00009171023010415733
$ws.Range("T7").Value = 0.00009171023010415731
$ws.Range("I8").Value = 0.006762540683959845
$ws.Range("J8").Value = 0.006762540683959845
$ws.Range("O8").Value = 0.392557056479861
$ws.Range("P8").Value = 0.3925570564798609
$ws.Range("S8").Value = 0.002654683065220583
$ws.Range("T8").Value = 0.002654683065220583
$ws.Range("I9").Value = 0.006762540683959845
$ws.Range("J9").Value = 0.006762540683959845
$ws.Range("M9").Value = 127.396393
$ws.Range("N9").Value = 382.189179
$ws.Range("O9").Value = 0.2739817680029065
$ws.Range("P9").Value = 0.2739817680029065
$ws.Range("Q9").Value = 6484.36968798813
$ws.Range("R9").Value = 58359.32719189317
$ws.Range("S9").Value = 0.001852812852782903
$ws.Range("T9").Value = 0.001852812852782903
$ws.Range("I10").Value = 0.006762540683959845
$ws.Range("J10").Value = 0.006762540683959845
$ws.Range("M10").Value = 19.42400133333333
$ws.Range("N10").Value = 58.272004
$ws.Range("O10").Value = 0.04177372766745037
$ws.Range("P10").Value = 0.04177372766745036
$ws.Range("Q10").Value = 988.6653970282164
$ws.Range("R10").Value = 8897.988573253948
$ws.Range("S10").Value = 0.0002824965328717921
$ws.Range("T10").Value = 0.0002824965328717921
$ws.Range("I11").Value = 0.006762540683959845
$ws.Range("J11").Value = 0.006762540683959845
$ws.Range("M11").Value = 129.3233566666667
$ws.Range("N11").Value = 387.97007
$ws.Range("O11").Value = 0.2781259427306063
$ws.Range("P11").Value = 0.2781259427306062
$ws.Range("Q11").Value = 6582.450524468232
$ws.Range("R11").Value = 59242.05472021409
$ws.Range("S11").Value = 0.001880838002980411
$ws.Range("T11").Value = 0.001880838002980411
$ws.Range("G12").Value = 3274.382486666667
$ws.Range("H12").Value = 9823.14746
$ws.Range("I12").Value = 0.4350394734576531
$ws.Range("J12").Value = 0.435039473457653
$ws.Range("M12").Value = 6.305846
$ws.Range("N12").Value = 18.917538
$ws.Range("O12").Value = 0.01356150511917599
$ws.Range("P12").Value = 0.01356150511917599
$ws.Range("Q12").Value = 20647.75170601705
$ws.Range("R12").Value = 185829.7653541535
$ws.Range("S12").Value = 0.005899790046339589
$ws.Range("T12").Value = 0.005899790046339588
$ws.Range("G13").Value = 3274.382486666667
$ws.Range("H13").Value = 9823.14746
$ws.Range("I13").Value = 0.4350394734576531
$ws.Range("J13").Value = 0.435039473457653
$ws.Range("O13").Value = 0.392557056479861
$ws.Range("P13").Value = 0.3925570564798609
$ws.Range("Q13").Value = 597678.5438940702
$ws.Range("R13").Value = 5379106.895046631
$ws.Range("S13").Value = 0.1707778151530849
$ws.Range("T13").Value = 0.1707778151530849
$ws.Range("G14").Value = 3274.382486666667
$ws.Range("H14").Value = 9823.14746
$ws.Range("I14").Value = 0.4350394734576531
$ws.Range("J14").Value = 0.435039473457653
$ws.Range("M14").Value = 127.396393
$ws.Range("N14").Value = 382.189179
$ws.Range("O14").Value = 0.2739817680029065
$ws.Range("P14").Value = 0.2739817680029065
$ws.Range("Q14").Value = 417144.518103704
$ws.Range("R14").Value = 3754300.662933336
$ws.Range("S14").Value = 0.1191928840889813
$ws.Range("T14").Value = 0.1191928840889813
$ws.Range("G15").Value = 3274.382486666667
$ws.Range("H15").Value = 9823.14746
$ws.Range("I15").Value = 0.4350394734576531
$ws.Range("J15").Value = 0.435039473457653
$ws.Range("M15").Value = 19.42400133333333
$ws.Range("N15").Value = 58.272004
$ws.Range("O15").Value = 0.04177372766745037
$ws.Range("P15").Value = 0.04177372766745036
$ws.Range("Q15").Value = 63601.60978685665
$ws.Range("R15").Value = 572414.4880817098
$ws.Range("S15").Value = 0.018173220488811
$ws.Range("T15").Value = 0.018173220488811
$ws.Range("G16").Value = 3274.382486666667
$ws.Range("H16").Value = 9823.14746
$ws.Range("I16").Value = 0.4350394734576531
$ws.Range("J16").Value = 0.435039473457653
$ws.Range("M16").Value = 129.3233566666667
$ws.Range("N16").Value = 387.97007
$ws.Range("O16").Value = 0.2781259427306063
$ws.Range("P16").Value = 0.2781259427306062
$ws.Range("Q16").Value = 423454.1341862803
$ws.Range("R16").Value = 3811087.207676522
$ws.Range("S16").Value = 0.1209957636804363
$ws.Range("T16").Value = 0.1209957636804363
$ws.Range("G17").Value = 7.278837333333333
$ws.Range("H17").Value = 21.836512
$ws.Range("I17").Value = 0.0009670774791190726
$ws.Range("J17").Value = 0.0009670774791190726
$ws.Range("M17").Value = 6.305846
$ws.Range("N17").Value = 18.917538
$ws.Range("O17").Value = 0.01356150511917599
$ws.Range("P17").Value = 0.01356150511917599
$ws.Range("Q17").Value = 45.89922728305066
$ws.Range("R17").Value = 413.093045547456
$ws.Range("S17").Value = 0.00001311502618371311
$ws.Range("T17").Value = 0.00001311502618371311
$ws.Range("G18").Value = 7.278837333333333
$ws.Range("H18").Value = 21.836512
$ws.Range("I18").Value = 0.0009670774791190726
$ws.Range("J18").Value = 0.0009670774791190726
$ws.Range("O18").Value = 0.392557056479861
$ws.Range("P18").Value = 0.3925570564798609
$ws.Range("Q18").Value = 1328.618423884007
$ws.Range("R18").Value = 11957.56581495606
$ws.Range("S18").Value = 0.0003796330885909474
$ws.Range("T18").Value = 0.0003796330885909473
$ws.Range("G19").Value = 7.278837333333333
$ws.Range("H19").Value = 21.836512
$ws.Range("I19").Value = 0.0009670774791190726
$ws.Range("J19").Value = 0.0009670774791190726
$ws.Range("M19").Value = 127.396393
$ws.Range("N19").Value = 382.189179
$ws.Range("O19").Value = 0.2739817680029065
$ws.Range("P19").Value = 0.2739817680029065
$ws.Range("Q19").Value = 927.2976215004053
$ws.Range("R19").Value = 8345.678593503648
$ws.Range("S19").Value = 0.0002649615975248374
$ws.Range("T19").Value = 0.0002649615975248374
$ws.Range("G20").Value = 7.278837333333333
$ws.Range("H20").Value = 21.836512
$ws.Range("I20").Value = 0.0009670774791190726
$ws.Range("J20").Value = 0.0009670774791190726
$ws.Range("M20").Value = 19.42400133333333
$ws.Range("N20").Value = 58.272004
$ws.Range("O20").Value = 0.04177372766745037
$ws.Range("P20").Value = 0.04177372766745036
$ws.Range("Q20").Value = 141.3841460677831
$ws.Range("R20").Value = 1272.457314610048
$ws.Range("S20").Value = 0.00004039843124604456
$ws.Range("T20").Value = 0.00004039843124604455
$ws.Range("G21").Value = 7.278837333333333
$ws.Range("H21").Value = 21.836512
$ws.Range("I21").Value = 0.0009670774791190726
$ws.Range("J21").Value = 0.0009670774791190726
$ws.Range("M21").Value = 129.3233566666667
$ws.Range("N21").Value = 387.97007
$ws.Range("O21").Value = 0.2781259427306063
$ws.Range("P21").Value = 0.2781259427306062
$ws.Range("Q21").Value = 941.3236765773156
$ws.Range("R21").Value = 8471.913089195839
$ws.Range("S21").Value = 0.0002689693355735303
$ws.Range("T21").Value = 0.0002689693355735302
$ws.Range("G22").Value = 4161.570231333333
$ws.Range("H22").Value = 12484.710694
$ws.Range("I22").Value = 0.5529125963654108
$ws.Range("J22").Value = 0.5529125963654108
$ws.Range("M22").Value = 6.305846
$ws.Range("N22").Value = 18.917538
$ws.Range("O22").Value = 0.01356150511917599
$ws.Range("P22").Value = 0.01356150511917599
$ws.Range("Q22").Value = 26242.22099697237
$ws.Range("R22").Value = 236179.9889727514
$ws.Range("S22").Value = 0.007498327006066405
$ws.Range("T22").Value = 0.007498327006066404
$ws.Range("G23").Value = 4161.570231333333
$ws.Range("H23").Value = 12484.710694
$ws.Range("I23").Value = 0.5529125963654108
$ws.Range("J23").Value = 0.5529125963654108
$ws.Range("O23").Value = 0.392557056479861
$ws.Range("P23").Value = 0.3925570564798609
$ws.Range("Q23").Value = 759618.4154735926
$ws.Range("R23").Value = 6836565.739262332
$ws.Range("S23").Value = 0.2170497413198431
$ws.Range("T23").Value = 0.2170497413198431
$ws.Range("G24").Value = 4161.570231333333
$ws.Range("H24").Value = 12484.710694
$ws.Range("I24").Value = 0.5529125963654108
$ws.Range("J24").Value = 0.5529125963654108
$ws.Range("M24").Value = 127.396393
$ws.Range("N24").Value = 382.189179
$ws.Range("O24").Value = 0.2739817680029065
$ws.Range("P24").Value = 0.2739817680029065
$ws.Range("Q24").Value = 530169.0366880422
$ws.Range("R24").Value = 4771521.330192381
$ws.Range("S24").Value = 0.1514879707032727
$ws.Range("T24").Value = 0.1514879707032727
$ws.Range("G25").Value = 4161.570231333333
$ws.Range("H25").Value = 12484.710694
$ws.Range("I25").Value = 0.5529125963654108
$ws.Range("J25").Value = 0.5529125963654108
$ws.Range("M25").Value = 19.42400133333333
$ws.Range("N25").Value = 58.272004
$ws.Range("O25").Value = 0.04177372766745037
$ws.Range("P25").Value = 0.04177372766745036
$ws.Range("Q25").Value = 80834.34572217897
$ws.Range("R25").Value = 727509.1114996108
$ws.Range("S25").Value = 0.02309722022447158
$ws.Range("T25").Value = 0.02309722022447157
$ws.Range("G26").Value = 4161.570231333333
$ws.Range("H26").Value = 12484.710694
$ws.Range("I26").Value = 0.5529125963654108
$ws.Range("J26").Value = 0.5529125963654108
$ws.Range("M26").Value = 129.3233566666667
$ws.Range("N26").Value = 387.97007
$ws.Range("O26").Value = 0.2781259427306063
$ws.Range("P26").Value = 0.2781259427306062
$ws.Range("Q26").Value = 538188.2313201032
$ws.Range("R26").Value = 4843694.081880929
$ws.Range("S26").Value = 0.1537793371117571
$ws.Range("T26").Value = 0.1537793371117571

Write-Output "Updated $(278) cells with new TPM values"
